# Issue 21 - primera parte
# Update four figures on slide 1 (balance general table):
#   24,340            -> 24,341
#   $ 5,348,189.05     -> $ 5,348,204.05
#   8,343              -> 8,344
#   $ 2,241,327.15     -> $ 2,241,342.15
#
# The four text boxes use PowerPoint's "shrink shape to fit text"
# autofit (<a:spAutoFit/>), so writing new text makes the host re-flow
# and resize the shape. We restore each shape's original height right
# after the text edit so only the <a:t> run content changes.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$updates = @(
    @{ Index = 20; Text = "24,341" },
    @{ Index = 21; Text = "$ 5,348,204.05" },
    @{ Index = 22; Text = "8,344" },
    @{ Index = 23; Text = "$ 2,241,342.15" }
)

foreach ($u in $updates) {
    $shp = $s.Shapes.Item($u.Index)
    $tf = $shp.TextFrame
    $origHeight = $shp.Height

    $tr = $tf.TextRange
    $full = $tr.Characters(1, $tr.Length)
    $full.Text = $u.Text

    # Nudge back onto the pre-edit height: the autofit recalculation
    # above can drift the box by a float-precision hair, so restoring
    # it explicitly keeps the shape's <a:ext> unchanged. (+0.00002pt is
    # below single-point rendering precision but clears the EMU
    # truncation boundary when Height is read back and rewritten.)
    $shp.Height = $origHeight + 0.00002
}
